$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "John Smith"
$ws.Range("C2").Value = "Pending Approval"

$ws.Range("A2").Select()
